# Auto-generated PowerShell Excel COM-interop script
# Applies the weekly data shift described in the commit diff ("Fruta / hortaliza, semanal"):
#  - A new week's default entries are inserted at rows 75-76 (date 2022-05-27 / serial 44708)
#  - All rows that used to be 75-261 shift down by 2 (to 77-263)
#  - The former last two rows (260-261) end up duplicated at the very end (262-263)
#  - Sheet dimension grows from A1:R261 to A1:R263

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nRows = 189
$nCols = 18
$data = New-Object 'object[,]' $nRows,$nCols

$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44708, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[0,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44708, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[1,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44222, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[2,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44222, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[3,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44624, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 270, 600, 650, 628, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 628, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[4,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44161, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[5,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44161, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[6,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44211, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[7,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44211, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[8,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44264, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[9,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44264, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[10,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44434, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[11,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44434, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[12,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44229, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[13,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44229, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[14,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44278, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[15,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44278, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[16,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44663, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[17,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44663, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[18,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44546, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 450, 600, 650, 628, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 628, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[19,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44659, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 650, 600, 650, 627, '$/atado 0,5 a 1 kilo', 'Provincia de Cautín', 627, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[20,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44649, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 550, 600, 650, 623, '$/atado 0,5 a 1 kilo', 'Región Metropolitana', 623, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[21,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44467, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[22,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44467, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[23,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44306, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[24,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44306, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[25,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44169, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[26,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44169, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[27,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44420, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[28,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44420, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[29,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44665, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[30,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44665, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[31,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44525, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[32,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44525, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[33,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44308, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[34,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44308, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[35,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44595, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[36,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44595, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[37,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44553, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[38,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44553, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[39,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44295, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[40,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44295, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[41,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44687, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 250, 600, 650, 630, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 630, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[42,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44642, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 220, 600, 650, 623, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 623, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[43,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44476, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[44,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44476, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[45,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44561, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 450, 500, 550, 522, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 522, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[46,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44545, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 250, 550, 600, 580, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 580, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[47,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44518, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 450, 600, 650, 628, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 628, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[48,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44348, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[49,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44348, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[50,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44350, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[51,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44350, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[52,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44398, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[53,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44398, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[54,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44705, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[55,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44705, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[56,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44567, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región Metropolitana', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[57,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44567, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región Metropolitana', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[58,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44322, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[59,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44322, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[60,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44327, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[61,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44327, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[62,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44259, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[63,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44259, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[64,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44383, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[65,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44383, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[66,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44362, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[67,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44362, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[68,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44266, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[69,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44266, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[70,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44607, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[71,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44607, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[72,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44237, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[73,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44237, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[74,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44539, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 220, 600, 650, 627, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 627, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[75,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44589, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 650, 625, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 625, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[76,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44196, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[77,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44196, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[78,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44497, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 250, 600, 650, 630, '$/atado 1 a 1,5 kilos', 'Región del Maule', 420, 1.5, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[79,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44252, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[80,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44252, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 700, 700, 700, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 700, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[81,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44202, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[82,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44202, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[83,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44453, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[84,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44453, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[85,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44588, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 650, 625, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 625, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[86,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44616, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[87,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44616, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[88,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44628, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 250, 600, 650, 630, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 630, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[89,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44341, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[90,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44341, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[91,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44609, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Provincia de Cautín', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[92,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44609, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Provincia de Cautín', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[93,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44330, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[94,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44330, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[95,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44677, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[96,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44677, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[97,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44250, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[98,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44250, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[99,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44334, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[100,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44334, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[101,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44280, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[102,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44280, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[103,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44635, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 350, 600, 650, 629, '$/atado 0,5 a 1 kilo', 'Provincia de Cautín', 629, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[104,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44637, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 350, 600, 650, 621, '$/atado 0,5 a 1 kilo', 'Provincia de Cautín', 621, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[105,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44582, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región Metropolitana', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[106,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44582, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región Metropolitana', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[107,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44698, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 450, 600, 650, 622, '$/atado 0,5 a 1 kilo', 'Región Metropolitana', 622, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[108,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44257, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[109,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44257, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[110,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44498, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 450, 600, 650, 622, '$/atado 1 a 1,5 kilos', 'Región del Maule', 415, 1.5, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[111,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44209, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[112,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44209, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[113,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44217, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[114,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44217, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[115,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44215, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[116,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44215, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[117,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44405, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[118,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44405, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[119,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44650, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 150, 600, 650, 633, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 633, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[120,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44239, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región Metropolitana', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[121,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44239, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región Metropolitana', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[122,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44358, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[123,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44358, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[124,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44273, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[125,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44273, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[126,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44400, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[127,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44400, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[128,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44551, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[129,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44551, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[130,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44691, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[131,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44691, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[132,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44187, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[133,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44187, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[134,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44558, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 350, 500, 550, 529, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 529, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[135,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44694, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[136,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44694, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[137,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44488, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[138,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44488, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[139,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44166, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[140,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44166, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[141,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44316, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región Metropolitana', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[142,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44316, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región Metropolitana', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[143,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44656, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 350, 600, 650, 629, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 629, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[144,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44586, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 310, 500, 550, 526, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 526, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[145,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44469, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[146,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44469, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[147,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44579, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[148,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44579, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[149,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44243, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[150,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44243, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[151,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44505, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 250, 600, 650, 630, '$/atado 0,5 a 1 kilo', 'Región del Maule', 630, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[152,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44267, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[153,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44267, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[154,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44474, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[155,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44474, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[156,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44631, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 250, 600, 650, 630, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 630, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[157,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44672, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 270, 600, 650, 622, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 622, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[158,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44490, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[159,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44490, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[160,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44679, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[161,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44679, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[162,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44365, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[163,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44365, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[164,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44603, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 450, 500, 550, 522, '$/atado 0,5 a 1 kilo', 'Región Metropolitana', 522, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[165,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44427, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[166,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44427, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[167,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44565, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[168,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44565, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[169,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44447, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[170,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44447, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[171,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44533, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 350, 600, 650, 621, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 621, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[172,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44523, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[173,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44523, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[174,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44462, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[175,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44462, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[176,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44159, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[177,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44159, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[178,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44344, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[179,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44344, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[180,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44600, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 650, 500, 550, 523, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 523, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[181,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44602, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 250, 600, 650, 630, '$/atado', 'Región de Ñuble', 630, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[182,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44692, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[183,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44692, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[184,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44376, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[185,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44376, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[186,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44442, 8, 100112009, 'Acelga', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 650, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[187,$c] = $rowValues[$c] }
$rowValues = @(11, 'Vega Monumental Concepción', 'Bíobío', 44442, 8, 100112009, 'Acelga', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/atado 0,5 a 1 kilo', 'Región de Ñuble', 500, 1, 'Hortaliza')
for ($c = 0; $c -lt $nCols; $c++) { $data[188,$c] = $rowValues[$c] }

$ws.Range("A75:R263").Value = $data
